$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price / 1h-volume figures.
# Price cells whose new text still parses as a plain number (e.g. "1.00",
# "5.80") need an explicit Text number format first, otherwise Excel
# auto-converts the literal ("5.80" -> 5.8) and the formatting is lost.
$ws.Cells.Item(2, 4).Value = '36.411.63'
$ws.Cells.Item(2, 5).Value = '  -1.51%  '
$ws.Cells.Item(3, 4).Value = '2.038.31'
$ws.Cells.Item(3, 5).Value = '  +2.52%  '
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '230.77'
$ws.Cells.Item(5, 5).Value = '  -11.92%  '
$ws.Cells.Item(6, 5).Value = '  -1.23%  '
$ws.Cells.Item(7, 5).Value = '  -0.13%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '54.88'
$ws.Cells.Item(8, 5).Value = '  -1.19%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.369'
$ws.Cells.Item(9, 5).Value = '  -0.80%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '56.93'
$ws.Cells.Item(10, 5).Value = '  +1.49%  '
$ws.Cells.Item(11, 5).Value = '  -1.69%  '
$ws.Cells.Item(12, 5).Value = '  +0.14%  '
$ws.Cells.Item(13, 4).Value = '2.338.25'
$ws.Cells.Item(13, 5).Value = '  +2.52%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '14.29'
$ws.Cells.Item(14, 5).Value = '  +1.65%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '19.93'
$ws.Cells.Item(15, 5).Value = '  -8.80%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.756'
$ws.Cells.Item(16, 5).Value = '  -1.07%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '5.08'
$ws.Cells.Item(17, 5).Value = '  -1.25%  '
$ws.Cells.Item(18, 4).Value = '2.065.43'
$ws.Cells.Item(18, 5).Value = '  +3.94%  '
$ws.Cells.Item(19, 4).Value = '36.448.11'
$ws.Cells.Item(19, 5).Value = '  -1.04%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '5.80'
$ws.Cells.Item(20, 5).Value = '  +15.18%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '67.38'
$ws.Cells.Item(21, 5).Value = '  -2.93%  '
$ws.Cells.Item(22, 5).Value = '  -3.42%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '220.03'
$ws.Cells.Item(23, 5).Value = '  -5.54%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '1.00'
$ws.Cells.Item(24, 5).Value = '  +0.02%  '
$ws.Cells.Item(25, 5).Value = '  +1.46%  '
$ws.Cells.Item(26, 5).Value = '  -8.09%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '162.68'
$ws.Cells.Item(27, 5).Value = '  -1.39%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '8.68'
$ws.Cells.Item(28, 5).Value = '  -1.35%  '
$ws.Cells.Item(29, 5).Value = '  -0.41%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '18.89'
$ws.Cells.Item(30, 5).Value = '  -1.85%  '
$ws.Cells.Item(31, 5).Value = '  +3.87%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.116'
$ws.Cells.Item(32, 5).Value = '  -1.46%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.35'
$ws.Cells.Item(33, 5).Value = '  -3.37%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0600'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.48'
$ws.Cells.Item(35, 5).Value = '  +3.73%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '4.24'
$ws.Cells.Item(36, 5).Value = '  -1.81%  '
$ws.Cells.Item(37, 5).Value = '  -0.10%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.76'
$ws.Cells.Item(38, 5).Value = '  -2.51%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.77'
$ws.Cells.Item(39, 5).Value = '  +8.37%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.22'
$ws.Cells.Item(40, 5).Value = '  -5.66%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '4.56'
$ws.Cells.Item(41, 5).Value = '  +47.31%  '
$ws.Cells.Item(42, 5).Value = '  -3.93%  '
$ws.Cells.Item(43, 4).Value = '1.479.51'
$ws.Cells.Item(43, 5).Value = '  +3.09%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0937'
$ws.Cells.Item(44, 5).Value = '  +3.19%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '93.20'
$ws.Cells.Item(45, 5).Value = '  +5.22%  '
$ws.Cells.Item(46, 5).Value = '  -1.56%  '
$ws.Cells.Item(47, 5).Value = '  -4.98%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '15.47'
$ws.Cells.Item(48, 5).Value = '  +0.62%  '
$ws.Cells.Item(49, 5).Value = '  -1.76%  '
$ws.Cells.Item(50, 5).Value = '  -0.33%  '
$ws.Cells.Item(51, 5).Value = '  +2.28%  '
